$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest reporting period (column D) which shifts all subsequent
# periods (E:M) one column to the left (D:L), matching how the source
# workbook drops the oldest quarter when a new one is published.
$ws.Columns.Item(4).Delete()

# The newly exposed last column (M) needs the same "wide" column width used
# by the other right-most columns of each 3-column group (stored width 29,
# which the Excel object model reports/accepts as ColumnWidth 28.17).
$ws.Columns.Item(13).ColumnWidth = 28.17

# --- Row 8: financial period headers ---
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Range("I9").Value = "1402-02-30 (9)"
$ws.Range("M9").Value = "1402-02-30"

# --- Financial data rows: fix up cells whose recomputed value differs from
# a pure shift, and populate the newly added period column (M). ---
# Row 11: Sales
$ws.Range("M11").Value = 82896423
# Row 12: Cost of goods sold
$ws.Range("M12").Value = -60690973
# Row 13: Gross profit (loss)
$ws.Range("M13").Value = 22205450
# Row 14: General, administrative and organizational expenses
$ws.Range("M14").Value = -9429283
# Row 15: Receivables impairment expense
$ws.Range("M15").Value = 0
# Row 16: Net other operating income (expenses)
$ws.Range("M16").Value = 130375
# Row 17: Operating profit (loss)
$ws.Range("M17").Value = 12906542
# Row 18: Financial expenses
$ws.Range("M18").Value = -7800501
# Row 19: Net other non-operating income and expenses
$ws.Range("M19").Value = 27400940
# Row 20: Net profit (loss) from continuing operations before tax
$ws.Range("M20").Value = 32506981
# Row 21: Tax
$ws.Range("I21").Value = -4292063
$ws.Range("M21").Value = -747736
# Row 22: Net profit (loss) from continuing operations
$ws.Range("I22").Value = 20648381
$ws.Range("M22").Value = 31759245
# Row 23: Profit (loss) from discontinued operations, net of tax
$ws.Range("M23").Value = 0
# Row 24: Net profit (loss)
$ws.Range("I24").Value = 20648381
$ws.Range("M24").Value = 31759245
# Row 25: EPS after tax
$ws.Range("I25").Value = 21735
$ws.Range("M25").Value = 1026
# Row 26: Capital
$ws.Range("I26").Value = 950000
$ws.Range("M26").Value = 30950000
# Row 27: EPS based on latest capital
$ws.Range("M27").Value = 1026
